# Admin transaction data verify updated
# Append three new transaction ID rows (B5:B7) to the "Transactions" sheet,
# written as text so they land in the shared-string table (not as numbers).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Transactions")

$newValues = @("3373363968", "3367494400", "3315788544")
$startRow = 5

for ($i = 0; $i -lt $newValues.Length; $i++) {
    $row = $startRow + $i
    $cell = $ws.Range("B$row")

    # Format as Text first so Excel stores the numeric-looking string as a
    # string (shared string) rather than silently coercing it to a number.
    $cell.NumberFormat = "@"
    $cell.Value = $newValues[$i]
    # Restore the default "Normal" cell style so no new style/number-format
    # record is left attached to the cell (matches the original formatting).
    $cell.Style = "Normal"
}
